# DeveloperGuide: Updated diagram & description for Logic, Model, Storage &
# Versioned Tasketch to suit our project.
#
# The only content change on the (single) slide is the class-box label that
# used to read "AddressBook" (first line of a two-line "AddressBook / Parser"
# label) -- it is renamed to "TaskBook" so the diagram matches the renamed
# project (AddressBook-Level3 -> the team's "Tasketch"/task-tracking app).
# All other text on the slide, shape formatting, fills, lines, etc. must stay
# untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

$oldLabel = "AddressBook"
$newLabel = "TaskBook"

$target = $null
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tf = $shp.TextFrame
        if ($tf.HasText) {
            $txt = $tf.TextRange.Text
            if ($txt.StartsWith($oldLabel)) {
                $target = $shp
                break
            }
        }
    }
}

if ($target -ne $null) {
    # Replace just the "AddressBook" run's characters so every other
    # character run (and its run-level formatting: size, color, bold,
    # italics, etc.) in the shape -- including the second "Parser" line --
    # is left completely untouched.
    $tr = $target.TextFrame.TextRange
    $run = $tr.Characters(1, $oldLabel.Length)
    $run.Text = $newLabel
}
